$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Change in inventories
$ws.Range("B6").Value = -4181000.0
$ws.Range("C6").Value = 10922000.0
$ws.Range("D6").Value = -5558000.0
$ws.Range("E6").Value = -10580000.0
$ws.Range("F6").Value = -14129000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 468000000.0
$ws.Range("C8").Value = 471000000.0
$ws.Range("D8").Value = 366914000.0
$ws.Range("E8").Value = 260041000.0
$ws.Range("F8").Value = 177234000.0
